$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.414553642272949
$ws.Range("B1").Value = 3.511605501174927
$ws.Range("C1").Value = 2.773596525192261
$ws.Range("D1").Value = 2.27135443687439
$ws.Range("E1").Value = 1.541318893432617
